$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-22 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-23 Thursday", 2) | Out-Null
$d.Content.Find.Execute("84×12=1008", $true, $false, $false, $false, $false, $true, 1, $false, "66×42=2772", 2) | Out-Null
$d.Content.Find.Execute("14×20=280", $true, $false, $false, $false, $false, $true, 1, $false, "87×97=8439", 2) | Out-Null
$d.Content.Find.Execute("57×56=3192", $true, $false, $false, $false, $false, $true, 1, $false, "84×34=2856", 2) | Out-Null
$d.Content.Find.Execute("89×28=2492", $true, $false, $false, $false, $false, $true, 1, $false, "26×43=1118", 2) | Out-Null
$d.Content.Find.Execute("43×43=1849", $true, $false, $false, $false, $false, $true, 1, $false, "59×33=1947", 2) | Out-Null
$d.Content.Find.Execute("64×88=5632", $true, $false, $false, $false, $false, $true, 1, $false, "27×96=2592", 2) | Out-Null
$d.Content.Find.Execute("53×50=2650", $true, $false, $false, $false, $false, $true, 1, $false, "95×31=2945", 2) | Out-Null
$d.Content.Find.Execute("70×35=2450", $true, $false, $false, $false, $false, $true, 1, $false, "90×11=990", 2) | Out-Null
$d.Content.Find.Execute("36×77=2772", $true, $false, $false, $false, $false, $true, 1, $false, "88×13=1144", 2) | Out-Null
$d.Content.Find.Execute("38×37=1406", $true, $false, $false, $false, $false, $true, 1, $false, "68×41=2788", 2) | Out-Null
$d.Content.Find.Execute("52×97=5044", $true, $false, $false, $false, $false, $true, 1, $false, "81×62=5022", 2) | Out-Null
$d.Content.Find.Execute("35×65=2275", $true, $false, $false, $false, $false, $true, 1, $false, "80×84=6720", 2) | Out-Null
$d.Content.Find.Execute("43×70=3010", $true, $false, $false, $false, $false, $true, 1, $false, "82×99=8118", 2) | Out-Null
$d.Content.Find.Execute("58×99=5742", $true, $false, $false, $false, $false, $true, 1, $false, "94×55=5170", 2) | Out-Null
$d.Content.Find.Execute("35×50=1750", $true, $false, $false, $false, $false, $true, 1, $false, "75×64=4800", 2) | Out-Null
$d.Content.Find.Execute("40×41=1640", $true, $false, $false, $false, $false, $true, 1, $false, "16×78=1248", 2) | Out-Null
$d.Content.Find.Execute("60×68=4080", $true, $false, $false, $false, $false, $true, 1, $false, "85×12=1020", 2) | Out-Null
$d.Content.Find.Execute("46×77=3542", $true, $false, $false, $false, $false, $true, 1, $false, "64×97=6208", 2) | Out-Null
$d.Content.Find.Execute("49×23=1127", $true, $false, $false, $false, $false, $true, 1, $false, "93×47=4371", 2) | Out-Null
$d.Content.Find.Execute("79×11=869", $true, $false, $false, $false, $false, $true, 1, $false, "89×44=3916", 2) | Out-Null
$d.Content.Find.Execute("74×36=2664", $true, $false, $false, $false, $false, $true, 1, $false, "74×94=6956", 2) | Out-Null
$d.Content.Find.Execute("45×38=1710", $true, $false, $false, $false, $false, $true, 1, $false, "32×90=2880", 2) | Out-Null
$d.Content.Find.Execute("48×17=816", $true, $false, $false, $false, $false, $true, 1, $false, "33×88=2904", 2) | Out-Null
$d.Content.Find.Execute("20×18=360", $true, $false, $false, $false, $false, $true, 1, $false, "55×26=1430", 2) | Out-Null
$d.Content.Find.Execute("42×13=546", $true, $false, $false, $false, $false, $true, 1, $false, "36×82=2952", 2) | Out-Null
